$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before the existing "Late" column (old N, now becomes O)
$ws.Columns.Item(14).Insert()
$ws.Columns.Item(14).ColumnWidth = 10.17

# Make "Repayment schedule" the active tab/sheet and set its selection
$ws.Activate()
$ws.Range("K13").Select()
